# Updated cryptos list on Sat May 18 17:29:16 UTC 2024 with GitHub Actions
# Applies per-cell Price (D) and Volume(1h) (E) text updates to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Cell = 'D2'; Value = '66.947.45' },
    @{ Cell = 'E2'; Value = '  +0.02%  ' },
    @{ Cell = 'D3'; Value = '3.117.09' },
    @{ Cell = 'E3'; Value = '  +0.51%  ' },
    @{ Cell = 'E4'; Value = '  +0.05%  ' },
    @{ Cell = 'D5'; Value = '578.17' },
    @{ Cell = 'E5'; Value = '  -0.52%  ' },
    @{ Cell = 'D6'; Value = '171.46' },
    @{ Cell = 'E6'; Value = '  +1.84%  ' },
    @{ Cell = 'E7'; Value = '  +0.00%  ' },
    @{ Cell = 'D8'; Value = '3.114.85' },
    @{ Cell = 'E8'; Value = '  +0.59%  ' },
    @{ Cell = 'E9'; Value = '  -0.73%  ' },
    @{ Cell = 'D10'; Value = '6.48' },
    @{ Cell = 'E10'; Value = '  -2.65%  ' },
    @{ Cell = 'E11'; Value = '  -1.38%  ' },
    @{ Cell = 'D12'; Value = '0.483' },
    @{ Cell = 'E12'; Value = '  +0.10%  ' },
    @{ Cell = 'E13'; Value = '  -1.84%  ' },
    @{ Cell = 'D14'; Value = '37.25' },
    @{ Cell = 'E14'; Value = '  +1.07%  ' },
    @{ Cell = 'E15'; Value = '  -1.16%  ' },
    @{ Cell = 'D16'; Value = '3.635.38' },
    @{ Cell = 'E16'; Value = '  +0.55%  ' },
    @{ Cell = 'D17'; Value = '66.939.63' },
    @{ Cell = 'E18'; Value = '  -1.04%  ' },
    @{ Cell = 'D19'; Value = '3.116.96' },
    @{ Cell = 'E19'; Value = '  +0.52%  ' },
    @{ Cell = 'D20'; Value = '16.34' },
    @{ Cell = 'E20'; Value = '  -0.15%  ' },
    @{ Cell = 'D21'; Value = '477.11' },
    @{ Cell = 'E21'; Value = '  +1.45%  ' },
    @{ Cell = 'E22'; Value = '  -0.32%  ' },
    @{ Cell = 'D23'; Value = '7.97' },
    @{ Cell = 'E23'; Value = '  +5.49%  ' },
    @{ Cell = 'D24'; Value = '13.47' },
    @{ Cell = 'E24'; Value = '  +4.47%  ' },
    @{ Cell = 'E25'; Value = '  +0.80%  ' },
    @{ Cell = 'D26'; Value = '2.29' },
    @{ Cell = 'E26'; Value = '  -3.38%  ' },
    @{ Cell = 'D27'; Value = '10.10' },
    @{ Cell = 'E27'; Value = '  -0.82%  ' },
    @{ Cell = 'E28'; Value = '  +0.00%  ' },
    @{ Cell = 'E29'; Value = '  -2.10%  ' },
    @{ Cell = 'D30'; Value = '2.40' },
    @{ Cell = 'E30'; Value = '  -1.53%  ' },
    @{ Cell = 'E31'; Value = '  -0.24%  ' },
    @{ Cell = 'D32'; Value = '28.57' },
    @{ Cell = 'E32'; Value = '  +0.99%  ' },
    @{ Cell = 'E33'; Value = '  -0.06%  ' },
    @{ Cell = 'D34'; Value = '0.0₃0940' },
    @{ Cell = 'E34'; Value = '  -7.85%  ' },
    @{ Cell = 'E35'; Value = '  +0.09%  ' },
    @{ Cell = 'E36'; Value = '  -0.72%  ' },
    @{ Cell = 'D37'; Value = '0.973' },
    @{ Cell = 'D38'; Value = '47.01' },
    @{ Cell = 'E38'; Value = '  +0.26%  ' },
    @{ Cell = 'E39'; Value = '  -0.61%  ' },
    @{ Cell = 'D40'; Value = '50.06' },
    @{ Cell = 'E40'; Value = '  -0.44%  ' },
    @{ Cell = 'E41'; Value = '  -2.02%  ' },
    @{ Cell = 'E42'; Value = '  -0.99%  ' },
    @{ Cell = 'D43'; Value = '8.73' },
    @{ Cell = 'E43'; Value = '  -0.15%  ' },
    @{ Cell = 'D44'; Value = '2.851.15' },
    @{ Cell = 'E44'; Value = '  +3.25%  ' },
    @{ Cell = 'D45'; Value = '386.33' },
    @{ Cell = 'E45'; Value = '  -1.48%  ' },
    @{ Cell = 'E46'; Value = '  -8.32%  ' },
    @{ Cell = 'E47'; Value = '  -1.66%  ' },
    @{ Cell = 'D48'; Value = '136.16' },
    @{ Cell = 'E48'; Value = '  +0.95%  ' },
    @{ Cell = 'E49'; Value = '  +0.01%  ' },
    @{ Cell = 'D50'; Value = '24.90' },
    @{ Cell = 'E50'; Value = '  +0.34%  ' },
    @{ Cell = 'E51'; Value = '  -1.95%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force the assigned value to be stored as text, matching the source
    # data's inline-string cells (avoids Excel auto-coercing numeric-looking
    # strings such as "578.17" into numbers), then drop the temporary
    # number-format override so no stray cell style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
